$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4626922607421875
$ws.Range("E2").Value = 275.9337530369794
$ws.Range("F2").Value = 0.01040134582875035
$ws.Range("G2").Value = 0.008428035884090563
$ws.Range("H2").Value = 0.007738008719057779
$ws.Range("I2").Value = 0.007059611484734264
$ws.Range("J2").Value = 0.006821210903091132
$ws.Range("K2").Value = 0.006821210903091132
$ws.Range("L2").Value = 0.006418758235678772
$ws.Range("M2").Value = 0.006226746093999222
$ws.Range("N2").Value = 0.006213262495966666
$ws.Range("O2").Value = 0.006056658926514663
$ws.Range("P2").Value = 0.005833951560301759
$ws.Range("Q2").Value = 0.005833951560301759
$ws.Range("R2").Value = 0.005688151035449522
$ws.Range("S2").Value = 0.005636335713868044
$ws.Range("T2").Value = 0.005582386149234148
$ws.Range("U2").Value = 0.005541011573364975
$ws.Range("V2").Value = 0.005482114114696436
$ws.Range("W2").Value = 0.005440425200144848
$ws.Range("X2").Value = 0.005397674048609014
$ws.Range("Y2").Value = 0.005378825595262755
$ws.Range("C3").Value = 0.4530985355377197
$ws.Range("E3").Value = 262.8626510866325
$ws.Range("F3").Value = 0.01033874312311437
$ws.Range("G3").Value = 0.008500088864466711
$ws.Range("H3").Value = 0.007583552938507685
$ws.Range("I3").Value = 0.007257768139559727
$ws.Range("J3").Value = 0.006875546474686032
$ws.Range("K3").Value = 0.006655270305472172
$ws.Range("L3").Value = 0.00638282176669594
$ws.Range("M3").Value = 0.005927490665201376
$ws.Range("N3").Value = 0.005797455176214224
$ws.Range("O3").Value = 0.005745602826761158
$ws.Range("P3").Value = 0.005630781448283894
$ws.Range("Q3").Value = 0.005516013317684911
$ws.Range("R3").Value = 0.005417151072578647
$ws.Range("S3").Value = 0.005343764638332003
$ws.Range("T3").Value = 0.005317696231784359
$ws.Range("U3").Value = 0.005168841471454797
$ws.Range("V3").Value = 0.005168841471454797
$ws.Range("W3").Value = 0.005153457642781748
$ws.Range("X3").Value = 0.005134005563190174
$ws.Range("Y3").Value = 0.005124028286289131
$ws.Range("C4").Value = 0.3906357288360596
$ws.Range("E4").Value = 264.6186864901156
$ws.Range("F4").Value = 0.01047038990392477
$ws.Range("G4").Value = 0.008718343217483581
$ws.Range("H4").Value = 0.007736983885521965
$ws.Range("I4").Value = 0.007203006979048012
$ws.Range("J4").Value = 0.006982176386463974
$ws.Range("K4").Value = 0.006521126515472735
$ws.Range("L4").Value = 0.006267426779403299
$ws.Range("M4").Value = 0.006016747472192939
$ws.Range("N4").Value = 0.005880151418354016
$ws.Range("O4").Value = 0.005628025794764456
$ws.Range("P4").Value = 0.005552756952432648
$ws.Range("Q4").Value = 0.005509311953508329
$ws.Range("R4").Value = 0.005475881152192883
$ws.Range("S4").Value = 0.00543707239136625
$ws.Range("T4").Value = 0.005348516556342397
$ws.Range("U4").Value = 0.005315961966643704
$ws.Range("V4").Value = 0.005263303469299371
$ws.Range("W4").Value = 0.005177593693843964
$ws.Range("X4").Value = 0.005177593693843964
$ws.Range("Y4").Value = 0.005158258995908685
$ws.Range("C5").Value = 0.3593857288360596
$ws.Range("E5").Value = 256.8925975645307
$ws.Range("F5").Value = 0.01035150911078932
$ws.Range("G5").Value = 0.008673873058030284
$ws.Range("H5").Value = 0.008017025101872225
$ws.Range("I5").Value = 0.00718724823327683
$ws.Range("J5").Value = 0.006556711731160039
$ws.Range("K5").Value = 0.00647506170236483
$ws.Range("L5").Value = 0.006114808910002167
$ws.Range("M5").Value = 0.005897189574030352
$ws.Range("N5").Value = 0.005824898676250046
$ws.Range("O5").Value = 0.005592079379281949
$ws.Range("P5").Value = 0.005455307039384372
$ws.Range("Q5").Value = 0.005323865871543272
$ws.Range("R5").Value = 0.005225664481815406
$ws.Range("S5").Value = 0.005225664481815406
$ws.Range("T5").Value = 0.005183425730518822
$ws.Range("U5").Value = 0.005106470168867748
$ws.Range("V5").Value = 0.005091627630257076
$ws.Range("W5").Value = 0.005055263008845876
$ws.Range("X5").Value = 0.005025304779889503
$ws.Range("Y5").Value = 0.005007652973967458
$ws.Range("C6").Value = 0.3906159400939941
$ws.Range("E6").Value = 267.7499187496887
$ws.Range("F6").Value = 0.01031406849820599
$ws.Range("G6").Value = 0.00871924635319505
$ws.Range("H6").Value = 0.007653007787958506
$ws.Range("I6").Value = 0.007156190955423947
$ws.Range("J6").Value = 0.006734252341349758
$ws.Range("K6").Value = 0.006396797063456932
$ws.Range("L6").Value = 0.00625204549694566
$ws.Range("M6").Value = 0.006006521547083767
$ws.Range("N6").Value = 0.005773282397745099
$ws.Range("O6").Value = 0.005672263232356863
$ws.Range("P6").Value = 0.005598575024120759
$ws.Range("Q6").Value = 0.005545093658202197
$ws.Range("R6").Value = 0.005520360311779811
$ws.Range("S6").Value = 0.005456280298791846
$ws.Range("T6").Value = 0.005327616626460229
$ws.Range("U6").Value = 0.005327616626460229
$ws.Range("V6").Value = 0.005308654907050054
$ws.Range("W6").Value = 0.005287397233005817
$ws.Range("X6").Value = 0.005251100138647935
$ws.Range("Y6").Value = 0.005219296661787304
$ws.Range("C7").Value = 0.3749988079071045
$ws.Range("E7").Value = 266.9597803666729
$ws.Range("F7").Value = 0.01023719290449155
$ws.Range("G7").Value = 0.00857832451016282
$ws.Range("H7").Value = 0.008130251148314662
$ws.Range("I7").Value = 0.007207252445886939
$ws.Range("J7").Value = 0.006923596013504867
$ws.Range("K7").Value = 0.0067772327188038
$ws.Range("L7").Value = 0.006502048975950481
$ws.Range("M7").Value = 0.006157090176892234
$ws.Range("N7").Value = 0.005943183652129765
$ws.Range("O7").Value = 0.005783965117215661
$ws.Range("P7").Value = 0.005729620273471732
$ws.Range("Q7").Value = 0.005599045958942685
$ws.Range("R7").Value = 0.005425316105652722
$ws.Range("S7").Value = 0.005425218957624694
$ws.Range("T7").Value = 0.005329686864908693
$ws.Range("U7").Value = 0.005326343804999726
$ws.Range("V7").Value = 0.005298272164332605
$ws.Range("W7").Value = 0.005251192350215222
$ws.Range("X7").Value = 0.005232868555127321
$ws.Range("Y7").Value = 0.005203894354126177
$ws.Range("C8").Value = 0.3593747615814209
$ws.Range("E8").Value = 259.8245911551476
$ws.Range("F8").Value = 0.01031900034229313
$ws.Range("G8").Value = 0.008498085470277904
$ws.Range("H8").Value = 0.007519637791178748
$ws.Range("I8").Value = 0.007017384069639477
$ws.Range("J8").Value = 0.006390258430665003
$ws.Range("K8").Value = 0.006268325914412807
$ws.Range("L8").Value = 0.006030083442461903
$ws.Range("M8").Value = 0.005989265336944727
$ws.Range("N8").Value = 0.005849155181809057
$ws.Range("O8").Value = 0.005637049498475932
$ws.Range("P8").Value = 0.005423753720056552
$ws.Range("Q8").Value = 0.005423753720056552
$ws.Range("R8").Value = 0.005236274272748475
$ws.Range("S8").Value = 0.005147935543924322
$ws.Range("T8").Value = 0.005147935543924322
$ws.Range("U8").Value = 0.005147935543924322
$ws.Range("V8").Value = 0.005119721079399486
$ws.Range("W8").Value = 0.005077945315859441
$ws.Range("X8").Value = 0.005077945315859441
$ws.Range("Y8").Value = 0.005064806845129581
$ws.Range("C9").Value = 0.4649057388305664
$ws.Range("E9").Value = 267.8443427358488
$ws.Range("F9").Value = 0.01034234017830077
$ws.Range("G9").Value = 0.008850368464315896
$ws.Range("H9").Value = 0.007794286264484257
$ws.Range("I9").Value = 0.007396591152871948
$ws.Range("J9").Value = 0.006950027930222922
$ws.Range("K9").Value = 0.006670825049210926
$ws.Range("L9").Value = 0.00654345553773773
$ws.Range("M9").Value = 0.006418626887227455
$ws.Range("N9").Value = 0.006287225344243048
$ws.Range("O9").Value = 0.005805028149615073
$ws.Range("P9").Value = 0.005686742798919752
$ws.Range("Q9").Value = 0.005527506615362369
$ws.Range("R9").Value = 0.005527506615362369
$ws.Range("S9").Value = 0.005371803823856987
$ws.Range("T9").Value = 0.005361534212143419
$ws.Range("U9").Value = 0.005330029822072003
$ws.Range("V9").Value = 0.005310330785566698
$ws.Range("W9").Value = 0.005251835464597697
$ws.Range("X9").Value = 0.005245909061849652
$ws.Range("Y9").Value = 0.005221137285299196
$ws.Range("C10").Value = 0.3984847068786621
$ws.Range("E10").Value = 274.185441865875
$ws.Range("F10").Value = 0.01048508615881068
$ws.Range("G10").Value = 0.008712168689796233
$ws.Range("H10").Value = 0.007907265634125417
$ws.Range("I10").Value = 0.007363245376596948
$ws.Range("J10").Value = 0.007213270313554104
$ws.Range("K10").Value = 0.006874955267004289
$ws.Range("L10").Value = 0.006498918709401772
$ws.Range("M10").Value = 0.006122210936848193
$ws.Range("N10").Value = 0.006099590286100001
$ws.Range("O10").Value = 0.005891849388986479
$ws.Range("P10").Value = 0.005788733028199897
$ws.Range("Q10").Value = 0.005693417198700012
$ws.Range("R10").Value = 0.005635361133135072
$ws.Range("S10").Value = 0.005608439387591108
$ws.Range("T10").Value = 0.005584684360110435
$ws.Range("U10").Value = 0.00550783370184881
$ws.Range("V10").Value = 0.005473870583280489
$ws.Range("W10").Value = 0.005427884821505531
$ws.Range("X10").Value = 0.005393976394079621
$ws.Range("Y10").Value = 0.005344745455475144
$ws.Range("C11").Value = 0.3749897480010986
$ws.Range("E11").Value = 271.5231501078761
$ws.Range("F11").Value = 0.01024566124439139
$ws.Range("G11").Value = 0.008845960170942399
$ws.Range("H11").Value = 0.00712072162235441
$ws.Range("I11").Value = 0.00712072162235441
$ws.Range("J11").Value = 0.006789924672593769
$ws.Range("K11").Value = 0.006369899632067151
$ws.Range("L11").Value = 0.006361783083948128
$ws.Range("M11").Value = 0.006031188240395158
$ws.Range("N11").Value = 0.005893909539794811
$ws.Range("O11").Value = 0.005625657954778288
$ws.Range("P11").Value = 0.005625657954778288
$ws.Range("Q11").Value = 0.005556574985080888
$ws.Range("R11").Value = 0.005506589163614847
$ws.Range("S11").Value = 0.005499444570437847
$ws.Range("T11").Value = 0.005388962723291877
$ws.Range("U11").Value = 0.005377777034938407
$ws.Range("V11").Value = 0.005359040908180871
$ws.Range("W11").Value = 0.005329388010334644
$ws.Range("X11").Value = 0.005312849150142111
$ws.Range("Y11").Value = 0.005292848929978089
